$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the SAJE / GIAD BRIEF result values (columns B:E, rows 2:13)
$ws.Range("B2:E13").ClearContents()

# Update the selection to match the post-edit state
$ws.Range("B2:E14").Select()
